$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Column C ("Förändrad") holds a date serial number for every data row
# (rows 2 through 470). The whole column was bumped forward by one day
# (45178 -> 45179), so update every populated cell in that range.
$range = $ws.Range("C2:C470")
$range.Value = 45179
